$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7467.4287
$ws.Range("I28").Value = 1043.25
$ws.Range("J28").Value = 16033
$ws.Range("K28").Value = 1043.25
$ws.Range("L28").Value = 16033
$ws.Range("M28").Value = -558.25
$ws.Range("N28").Value = -17003
$ws.Range("H58").Value = 1929.875
$ws.Range("J58").Value = 2999.4
$ws.Range("L58").Value = 8998.200000000001
$ws.Range("N58").Value = -9298.200000000001
$ws.Range("H62").Value = 7650.875
$ws.Range("I62").Value = 4313.5
$ws.Range("J62").Value = 10988.25
$ws.Range("K62").Value = 4313.5
$ws.Range("L62").Value = 10988.25
$ws.Range("M62").Value = -3689.5
$ws.Range("N62").Value = -12236.25
$ws.Range("H65").Value = 7650.875
$ws.Range("I65").Value = 4313.5
$ws.Range("J65").Value = 10988.25
$ws.Range("K65").Value = 21567.5
$ws.Range("L65").Value = 54941.25
$ws.Range("M65").Value = -18447.5
$ws.Range("N65").Value = -61181.25
$ws.Range("H98").Value = 1352.6316
$ws.Range("I98").Value = 1043.75
$ws.Range("K98").Value = 1043.75
$ws.Range("M98").Value = 454.25
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 5000
$ws.Range("M106").Value = -4369
$ws.Range("H113").Value = 6777.5
$ws.Range("I113").Value = 8651.25
$ws.Range("J113").Value = 4903.75
$ws.Range("K113").Value = 8651.25
$ws.Range("L113").Value = 4903.75
$ws.Range("M113").Value = -5397.25
$ws.Range("N113").Value = -11411.75
$ws.Range("H116").Value = 8746.75
$ws.Range("J116").Value = 4299.5
$ws.Range("L116").Value = 4299.5
$ws.Range("N116").Value = -11183.5
$ws.Range("H122").Value = 1352.6316
$ws.Range("I122").Value = 1043.75
$ws.Range("K122").Value = 3131.25
$ws.Range("M122").Value = -681.25
$ws.Range("H132").Value = 6970.1
$ws.Range("J132").Value = 17133.334
$ws.Range("L132").Value = 51400.00199999999
$ws.Range("N132").Value = -56460.00199999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 5020177.5
$ws.Range("J24").Value = 5020177.5
$ws.Range("L24").Value = 5020177.5
$ws.Range("N24").Value = -5020925.5
$ws.Range("H32").Value = 4848.6113
$ws.Range("I32").Value = 4848.6113
$ws.Range("K32").Value = 4848.6113
$ws.Range("M32").Value = -4561.6113
$ws.Range("H96").Value = 3366855
$ws.Range("J96").Value = 3366855
$ws.Range("L96").Value = 3366855
$ws.Range("N96").Value = -3372347
$ws.Range("H100").Value = 5020177.5
$ws.Range("J100").Value = 5020177.5
$ws.Range("L100").Value = 5020177.5
$ws.Range("N100").Value = -5022341.5
$ws.Range("H110").Value = 4318.3335
$ws.Range("I110").Value = 3727.5
$ws.Range("J110").Value = 5500
$ws.Range("K110").Value = 3727.5
$ws.Range("L110").Value = 5500
$ws.Range("M110").Value = -1682.5
$ws.Range("N110").Value = -9590
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2805.1516
$ws.Range("I107").Value = 782.8
$ws.Range("K107").Value = 782.8
$ws.Range("M107").Value = 1137.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 918.2727
$ws.Range("I22").Value = 789.2222
$ws.Range("J22").Value = 1499
$ws.Range("K22").Value = 789.2222
$ws.Range("L22").Value = 1499
$ws.Range("M22").Value = -439.2222
$ws.Range("N22").Value = -2199
$ws.Range("H28").Value = 19500
$ws.Range("J28").Value = 19500
$ws.Range("L28").Value = 19500
$ws.Range("N28").Value = -19990
$ws.Range("H92").Value = 39999
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H122").Value = 868.86664
$ws.Range("I122").Value = 681
$ws.Range("K122").Value = 2043
$ws.Range("M122").Value = 407
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 260.2353
$ws.Range("J12").Value = 274.5
$ws.Range("L12").Value = 823.5
$ws.Range("N12").Value = -1169.5
$ws.Range("H14").Value = 140.5
$ws.Range("I14").Value = 140.5
$ws.Range("K14").Value = 421.5
$ws.Range("M14").Value = -248.5
$ws.Range("H109").Value = 2758.3333
$ws.Range("I109").Value = 1480
$ws.Range("K109").Value = 4440
$ws.Range("M109").Value = -3400
$ws.Range("H121").Value = 1427
$ws.Range("J121").Value = 2499.8572
$ws.Range("L121").Value = 7499.571599999999
$ws.Range("N121").Value = -10119.5716
$ws.Range("H131").Value = 1947.8572
$ws.Range("J131").Value = 1100
$ws.Range("L131").Value = 3300
$ws.Range("N131").Value = -13380
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 99008
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 99008
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 99008
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -100648
$ws.Range("H95").Value = 27672
$ws.Range("J95").Value = 27672
$ws.Range("L95").Value = 27672
$ws.Range("N95").Value = -33164
$ws.Range("H107").Value = 393.8889
$ws.Range("I107").Value = 363.57144
$ws.Range("K107").Value = 363.57144
$ws.Range("M107").Value = 1556.42856
$ws.Range("H134").Value = 117662.5
$ws.Range("J134").Value = 117662.5
$ws.Range("L134").Value = 352987.5
$ws.Range("N134").Value = -358057.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4664.25
$ws.Range("I7").Value = 3021.5454
$ws.Range("K7").Value = 3021.5454
$ws.Range("M7").Value = -2909.5454
$ws.Range("H16").Value = 385.44446
$ws.Range("I16").Value = 371.125
$ws.Range("K16").Value = 371.125
$ws.Range("M16").Value = -201.125
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 8000
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -7864
$ws.Range("N40").Value = -8272
$ws.Range("H46").Value = 5717.2104
$ws.Range("I46").Value = 3958.5557
$ws.Range("J46").Value = 7300
$ws.Range("K46").Value = 3958.5557
$ws.Range("L46").Value = 7300
$ws.Range("M46").Value = -3770.5557
$ws.Range("N46").Value = -7676
$ws.Range("H126").Value = 4664.25
$ws.Range("I126").Value = 3021.5454
$ws.Range("K126").Value = 9064.636200000001
$ws.Range("M126").Value = -6594.636200000001
$ws.Range("H132").Value = 4467.8184
$ws.Range("I132").Value = 3921.7778
$ws.Range("K132").Value = 11765.3334
$ws.Range("M132").Value = -9235.3334
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 8869.25
$ws.Range("J15").Value = 7737.5
$ws.Range("L15").Value = 7737.5
$ws.Range("N15").Value = -8313.5
$ws.Range("H69").Value = 60000
$ws.Range("J69").Value = 60000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61498
$ws.Range("H72").Value = 60000
$ws.Range("J72").Value = 60000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -187488
$ws.Range("H97").Value = 25572
$ws.Range("J97").Value = 25572
$ws.Range("L97").Value = 25572
$ws.Range("N97").Value = -27554
$ws.Range("H122").Value = 1749.5
$ws.Range("I122").Value = 2499
$ws.Range("K122").Value = 7497
$ws.Range("M122").Value = -5047
$ws.Range("H126").Value = 5954.364
$ws.Range("J126").Value = 6888.6665
$ws.Range("L126").Value = 20665.9995
$ws.Range("N126").Value = -25605.9995
